# Update timing trial results produced by the (now-fixed) algorithms.py runs.
# "Brute Force" and "Divide and Conquer" sheets get new per-trial timings;
# the "Summary" sheet recalculates automatically via its existing formulas.
$wb = $excel.ActiveWorkbook

$wsBF = $wb.Worksheets.Item("Brute Force")
$wsBF.Range("B2").Value = 0.000004053115844726561652967053
$wsBF.Range("C2").Value = 0.000001668930053710937923516474
$wsBF.Range("D2").Value = 0.000001192092895507812076483526
$wsBF.Range("F2").Value = 0.0000007152557373046875
$wsBF.Range("G2").Value = 0.00000095367431640625
$wsBF.Range("H2").Value = 0.000001192092895507812076483526
$wsBF.Range("J2").Value = 0.000001192092895507812076483526
$wsBF.Range("B3").Value = 0.030915021896362301218053048046
$wsBF.Range("C3").Value = 0.0305478572845458984375
$wsBF.Range("D3").Value = 0.030556917190551761281946951954
$wsBF.Range("E3").Value = 0.030209064483642581594446951954
$wsBF.Range("F3").Value = 0.030548810958862301218053048046
$wsBF.Range("G3").Value = 0.0300691127777099609375
$wsBF.Range("H3").Value = 0.0299241542816162109375
$wsBF.Range("I3").Value = 0.030415773391723629343053048046
$wsBF.Range("J3").Value = 0.030765771865844730031946951954
$wsBF.Range("K3").Value = 0.0305969715118408203125
$wsBF.Range("B4").Value = 0.122916698455810505241636576557
$wsBF.Range("C4").Value = 0.121798276901245103309712192186
$wsBF.Range("D4").Value = 0.123498916625976604133363423443
$wsBF.Range("E4").Value = 0.122605800628662095497212192186
$wsBF.Range("F4").Value = 0.122128963470458998252787807814
$wsBF.Range("G4").Value = 0.121542930603027302116636576557
$wsBF.Range("H4").Value = 0.122481107711792006065287807814
$wsBF.Range("I4").Value = 0.122097969055175795127787807814
$wsBF.Range("J4").Value = 0.120030879974365206619424384371
$wsBF.Range("K4").Value = 0.1220710277557373046875
$wsBF.Range("B5").Value = 0.2753736972808837890625
$wsBF.Range("C5").Value = 0.275441169738769475738848768742
$wsBF.Range("D5").Value = 0.274941205978393610198651231258
$wsBF.Range("E5").Value = 0.275262832641601618011151231258
$wsBF.Range("F5").Value = 0.2751941680908203125
$wsBF.Range("G5").Value = 0.2751648426055908203125
$wsBF.Range("H5").Value = 0.275336980819702092926348768742
$wsBF.Range("I5").Value = 0.2754991054534912109375
$wsBF.Range("J5").Value = 0.275842905044555719573651231258
$wsBF.Range("K5").Value = 0.2750437259674072265625
$wsBF.Range("B6").Value = 0.489665031433105524261151231258
$wsBF.Range("C6").Value = 0.4875490665435791015625
$wsBF.Range("D6").Value = 0.490314960479736272613848768742
$wsBF.Range("E6").Value = 0.490442037582397516448651231258
$wsBF.Range("F6").Value = 0.515902996063232421875
$wsBF.Range("G6").Value = 0.498270988464355524261151231258
$wsBF.Range("H6").Value = 0.495732307434081975738848768742
$wsBF.Range("I6").Value = 0.489017009735107421875
$wsBF.Range("J6").Value = 0.502807140350341796875
$wsBF.Range("K6").Value = 0.487269163131713922698651231258
$wsBF.Range("B7").Value = 0.7674758434295654296875
$wsBF.Range("C7").Value = 0.7678558826446533203125
$wsBF.Range("D7").Value = 0.7687680721282958984375
$wsBF.Range("E7").Value = 0.769464969635009765625
$wsBF.Range("F7").Value = 0.770394802093505859375
$wsBF.Range("G7").Value = 0.7616460323333740234375
$wsBF.Range("H7").Value = 0.76916980743408203125
$wsBF.Range("I7").Value = 0.770390033721923828125
$wsBF.Range("J7").Value = 0.7577068805694580078125
$wsBF.Range("K7").Value = 0.77945709228515625

$wsDC = $wb.Worksheets.Item("Divide and Conquer")
$wsDC.Range("B2").Value = 0.000032901763916015618223736422
$wsDC.Range("C2").Value = 0.000009059906005859375
$wsDC.Range("D2").Value = 0.000004053115844726561652967053
$wsDC.Range("E2").Value = 0.00000286102294921875
$wsDC.Range("F2").Value = 0.000003099441528320312076483526
$wsDC.Range("G2").Value = 0.0000019073486328125
$wsDC.Range("H2").Value = 0.0000019073486328125
$wsDC.Range("I2").Value = 0.000002145767211914062076483526
$wsDC.Range("J2").Value = 0.000001668930053710937923516474
$wsDC.Range("K2").Value = 0.000002145767211914062076483526
$wsDC.Range("B3").Value = 0.001395940780639648003819131006
$wsDC.Range("C3").Value = 0.001183032989501952908159565503
$wsDC.Range("D3").Value = 0.0012609958648681640625
$wsDC.Range("E3").Value = 0.001319169998168945095659565503
$wsDC.Range("F3").Value = 0.001292705535888672091840434497
$wsDC.Range("G3").Value = 0.001212120056152343966840434497
$wsDC.Range("H3").Value = 0.001235008239746093966840434497
$wsDC.Range("I3").Value = 0.001299858093261718966840434497
$wsDC.Range("J3").Value = 0.001294136047363281033159565503
$wsDC.Range("K3").Value = 0.001221895217895508029340434497
$wsDC.Range("B4").Value = 0.002658128738403319878819131006
$wsDC.Range("C4").Value = 0.002577781677246094183680868994
$wsDC.Range("D4").Value = 0.002562761306762694878819131006
$wsDC.Range("E4").Value = 0.0025961399078369140625
$wsDC.Range("F4").Value = 0.002532243728637694878819131006
$wsDC.Range("G4").Value = 0.0025379657745361328125
$wsDC.Range("H4").Value = 0.002564191818237305121180868994
$wsDC.Range("I4").Value = 0.002683162689208983941319131006
$wsDC.Range("J4").Value = 0.002639055252075194878819131006
$wsDC.Range("K4").Value = 0.002562046051025391058680868994
$wsDC.Range("B5").Value = 0.003829002380371094183680868994
$wsDC.Range("C5").Value = 0.003823757171630858941319131006
$wsDC.Range("D5").Value = 0.003704071044921875
$wsDC.Range("E5").Value = 0.003737688064575194878819131006
$wsDC.Range("F5").Value = 0.0036113262176513671875
$wsDC.Range("G5").Value = 0.003628015518188476996180868994
$wsDC.Range("H5").Value = 0.003867864608764648003819131006
$wsDC.Range("I5").Value = 0.003838062286376953125
$wsDC.Range("J5").Value = 0.0036571025848388671875
$wsDC.Range("K5").Value = 0.003650188446044921875
$wsDC.Range("B6").Value = 0.005367040634155272570138262012
$wsDC.Range("C6").Value = 0.0055510997772216796875
$wsDC.Range("D6").Value = 0.005207061767578125
$wsDC.Range("E6").Value = 0.0052440166473388671875
$wsDC.Range("F6").Value = 0.0053808689117431640625
$wsDC.Range("G6").Value = 0.00535106658935546875
$wsDC.Range("H6").Value = 0.00527286529541015625
$wsDC.Range("I6").Value = 0.005443096160888671875
$wsDC.Range("J6").Value = 0.0056436061859130859375
$wsDC.Range("K6").Value = 0.005290985107421875
$wsDC.Range("B7").Value = 0.00590801239013671875
$wsDC.Range("C7").Value = 0.0060908794403076171875
$wsDC.Range("D7").Value = 0.005795001983642578125
$wsDC.Range("E7").Value = 0.00586986541748046875
$wsDC.Range("F7").Value = 0.005669116973876953125
$wsDC.Range("G7").Value = 0.0058529376983642578125
$wsDC.Range("H7").Value = 0.005935907363891602429861737988
$wsDC.Range("I7").Value = 0.005820989608764647570138262012
$wsDC.Range("J7").Value = 0.0057380199432373046875
$wsDC.Range("K7").Value = 0.0057351589202880859375
